# Slide 10 ("Expressions versus Statements" / "Definition of CPRL" deck):
# in the body placeholder, the last paragraph's closing run
#   " function in Kotlin) and a procedure call is considered to be a statement."
# is split into three runs so that a comma is inserted right after "Kotlin)":
#   " function in " | "Kotlin), " | "and a procedure call is considered to be a statement."

$oldTail = " function in Kotlin) and a procedure call is considered to be a statement."
$part1   = " function in "
$part2   = "Kotlin), "
$part3   = "and a procedure call is considered to be a statement."

$p = $ppt.ActivePresentation

# Find the shape that actually contains the sentence we need to edit instead of
# hard-coding a shape index, so the script is resilient to shape-order changes.
$targetShape = $null
$targetSlide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text.IndexOf($oldTail) -ge 0) {
                $targetShape = $shape
                $targetSlide = $slide
                break
            }
        }
    }
    if ($targetShape -ne $null) {
        break
    }
}

$tr = $targetShape.TextFrame.TextRange

$fullText = $tr.Text
$charIndex = $fullText.IndexOf($oldTail)
$startPos = $charIndex + 1

# Replace the whole matched run's text in one shot (still a single run at this
# point) so the comma ends up in the text.
$whole = $tr.Characters($startPos, $oldTail.Length)
$whole.Text = $part1 + $part2 + $part3

# Now split that (now longer) stretch of text into three separate runs that
# match the three target runs.
$seg1 = $tr.Characters($startPos, $part1.Length)
$seg1.Text = $part1

$seg2 = $tr.Characters($startPos + $part1.Length, $part2.Length)
$seg2.Text = $part2

$seg3 = $tr.Characters($startPos + $part1.Length + $part2.Length, $part3.Length)
$seg3.Text = $part3
